# chore: update Sheets via scheduled runner
# Refresh cached market-board derived figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) for a handful of leve rows across the
# job sheets, picking up newer Universalis price snapshots.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 41668970
$ws.Range("I113").Value = 71430230
$ws.Range("J113").Value = 3195.2
$ws.Range("K113").Value = 71430230
$ws.Range("L113").Value = 3195.2
$ws.Range("M113").Value = -71426976
$ws.Range("N113").Value = -9703.200000000001

$ws.Range("H132").Value = 1747.5312
$ws.Range("I132").Value = 1359.0333
$ws.Range("K132").Value = 4077.0999
$ws.Range("M132").Value = -1547.0999

$ws.Range("H137").Value = 2679629
$ws.Range("I137").Value = 1163805.1
$ws.Range("J137").Value = 7693507.5
$ws.Range("K137").Value = 3491415.3
$ws.Range("L137").Value = 23080522.5
$ws.Range("M137").Value = -3488865.3
$ws.Range("N137").Value = -23085622.5

$ws.Range("H138").Value = 1921.0807
$ws.Range("I138").Value = 1137.2222
$ws.Range("J138").Value = 3006.423
$ws.Range("K138").Value = 3411.6666
$ws.Range("L138").Value = 9019.269
$ws.Range("M138").Value = 1728.3334
$ws.Range("N138").Value = -19299.269

$ws.Range("H141").Value = 2509.9714
$ws.Range("I141").Value = 2409.6
$ws.Range("J141").Value = 2585.25
$ws.Range("K141").Value = 7228.799999999999
$ws.Range("L141").Value = 7755.75
$ws.Range("M141").Value = -2048.799999999999
$ws.Range("N141").Value = -18115.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2399.8333
$ws.Range("I61").Value = 2399.8333
$ws.Range("K61").Value = 2399.8333
$ws.Range("M61").Value = -2187.8333

$ws.Range("H74").Value = 828.6
$ws.Range("I74").Value = 784.0678
$ws.Range("J74").Value = 1266.5
$ws.Range("K74").Value = 784.0678
$ws.Range("L74").Value = 1266.5
$ws.Range("M74").Value = 89.93219999999997
$ws.Range("N74").Value = -3014.5

$ws.Range("H77").Value = 828.6
$ws.Range("I77").Value = 784.0678
$ws.Range("J77").Value = 1266.5
$ws.Range("K77").Value = 3920.339
$ws.Range("L77").Value = 6332.5
$ws.Range("M77").Value = 447.6610000000001
$ws.Range("N77").Value = -15068.5

$ws.Range("H110").Value = 2581.9167
$ws.Range("I110").Value = 4298.727
$ws.Range("J110").Value = 1129.2307
$ws.Range("K110").Value = 4298.727
$ws.Range("L110").Value = 1129.2307
$ws.Range("M110").Value = -2253.727
$ws.Range("N110").Value = -5219.2307

$ws.Range("H132").Value = 159141.44
$ws.Range("I132").Value = 174121.6
$ws.Range("K132").Value = 522364.8
$ws.Range("M132").Value = -519834.8

$ws.Range("H136").Value = 2399.8333
$ws.Range("I136").Value = 2399.8333
$ws.Range("K136").Value = 7199.499899999999
$ws.Range("M136").Value = -4649.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 35780
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 35780
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 35780
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -36978

$ws.Range("H86").Value = 2253.1428
$ws.Range("I86").Value = 2300.889
$ws.Range("J86").Value = 1966.6666
$ws.Range("K86").Value = 2300.889
$ws.Range("L86").Value = 1966.6666
$ws.Range("M86").Value = -1177.889
$ws.Range("N86").Value = -4212.6666

$ws.Range("H89").Value = 2253.1428
$ws.Range("I89").Value = 2300.889
$ws.Range("J89").Value = 1966.6666
$ws.Range("K89").Value = 11504.445
$ws.Range("L89").Value = 9833.333000000001
$ws.Range("M89").Value = -5888.445
$ws.Range("N89").Value = -21065.333

$ws.Range("H134").Value = 167705.44
$ws.Range("I134").Value = 167705.44
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 503116.32
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -500581.32
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1810.4681
$ws.Range("I31").Value = 1429.8529
$ws.Range("J31").Value = 2805.923
$ws.Range("K31").Value = 1429.8529
$ws.Range("L31").Value = 2805.923
$ws.Range("M31").Value = -1134.8529
$ws.Range("N31").Value = -3395.923

$ws.Range("H34").Value = 1810.4681
$ws.Range("I34").Value = 1429.8529
$ws.Range("J34").Value = 2805.923
$ws.Range("K34").Value = 1429.8529
$ws.Range("L34").Value = 2805.923
$ws.Range("M34").Value = -1227.8529
$ws.Range("N34").Value = -3209.923

$ws.Range("H58").Value = 1920.1562
$ws.Range("I58").Value = 1949.8387
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 1949.8387
$ws.Range("L58").Value = 1000
$ws.Range("M58").Value = -1746.8387
$ws.Range("N58").Value = -1406

$ws.Range("H92").Value = 32233.334
$ws.Range("J92").Value = 32233.334
$ws.Range("L92").Value = 32233.334
$ws.Range("N92").Value = -37225.334

$ws.Range("H99").Value = 1740
$ws.Range("I99").Value = 1740
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1740
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -242
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 1740
$ws.Range("I126").Value = 1740
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5220
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2750
$ws.Range("N126").ClearContents()

$ws.Range("H134").Value = 6029.9473
$ws.Range("I134").Value = 6816.7095
$ws.Range("J134").Value = 2545.7144
$ws.Range("K134").Value = 20450.1285
$ws.Range("L134").Value = 7637.1432
$ws.Range("M134").Value = -17915.1285
$ws.Range("N134").Value = -12707.1432

$ws.Range("H136").Value = 1920.1562
$ws.Range("I136").Value = 1949.8387
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 5849.5161
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -3299.5161
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 327.53333
$ws.Range("I14").Value = 327.53333
$ws.Range("K14").Value = 982.5999899999999
$ws.Range("M14").Value = -809.5999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 11442.5
$ws.Range("I126").Value = 2960.7693
$ws.Range("J126").Value = 18793.334
$ws.Range("K126").Value = 8882.3079
$ws.Range("L126").Value = 56380.00199999999
$ws.Range("M126").Value = -6412.3079
$ws.Range("N126").Value = -61320.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2518.7222
$ws.Range("I7").Value = 2240.2
$ws.Range("K7").Value = 2240.2
$ws.Range("M7").Value = -2128.2

$ws.Range("H40").Value = 2197.5
$ws.Range("I40").Value = 2196.6667
$ws.Range("K40").Value = 2196.6667
$ws.Range("M40").Value = -2060.6667

$ws.Range("H126").Value = 2518.7222
$ws.Range("I126").Value = 2240.2
$ws.Range("K126").Value = 6720.599999999999
$ws.Range("M126").Value = -4250.599999999999

$ws.Range("H136").Value = 2035.8572
$ws.Range("I136").Value = 1961.0435
$ws.Range("J136").Value = 2380
$ws.Range("K136").Value = 5883.1305
$ws.Range("L136").Value = 7140
$ws.Range("M136").Value = -3333.1305
$ws.Range("N136").Value = -12240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 36709.93
$ws.Range("I113").Value = 62864.875
$ws.Range("J113").Value = 1836.6666
$ws.Range("K113").Value = 188594.625
$ws.Range("L113").Value = 5509.9998
$ws.Range("M113").Value = -186424.625
$ws.Range("N113").Value = -9849.9998

$ws.Range("H136").Value = 1960.5714
$ws.Range("I136").Value = 2160.5
$ws.Range("J136").Value = 1320.8
$ws.Range("K136").Value = 6481.5
$ws.Range("L136").Value = 3962.4
$ws.Range("M136").Value = -3931.5
$ws.Range("N136").Value = -9062.4
